$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column A from numeric auto-conversion while we write text values
$ws.Range("A2:A26").NumberFormat = "@"

# Header row: becomes text labels, keeps style s=1
$ws.Cells.Item(1,1).Value = "Price"
$ws.Cells.Item(1,2).Value = "Name"
$ws.Cells.Item(1,3).Value = "Symbol"

# Data rows: Price, Name, Symbol
$ws.Cells.Item(2,1).Value = "23,788.94"
$ws.Cells.Item(2,2).Value = "Bitcoin USD"
$ws.Cells.Item(2,3).Value = "BTC-USD"
$ws.Cells.Item(3,1).Value = "1,769.95"
$ws.Cells.Item(3,2).Value = "Ethereum USD"
$ws.Cells.Item(3,3).Value = "ETH-USD"
$ws.Cells.Item(4,1).Value = "1.0002"
$ws.Cells.Item(4,2).Value = "Tether USD"
$ws.Cells.Item(4,3).Value = "USDT-USD"
$ws.Cells.Item(5,1).Value = "0.999945"
$ws.Cells.Item(5,2).Value = "USD Coin USD"
$ws.Cells.Item(5,3).Value = "USDC-USD"
$ws.Cells.Item(6,1).Value = "324.47"
$ws.Cells.Item(6,2).Value = "Binance Coin USD"
$ws.Cells.Item(6,3).Value = "BNB-USD"
$ws.Cells.Item(7,1).Value = "0.376197"
$ws.Cells.Item(7,2).Value = "XRP USD"
$ws.Cells.Item(7,3).Value = "XRP-USD"
$ws.Cells.Item(8,1).Value = "0.533140"
$ws.Cells.Item(8,2).Value = "Cardano USD"
$ws.Cells.Item(8,3).Value = "ADA-USD"
$ws.Cells.Item(9,1).Value = "0.999966"
$ws.Cells.Item(9,2).Value = "Binance USD USD"
$ws.Cells.Item(9,3).Value = "BUSD-USD"
$ws.Cells.Item(10,1).Value = "42.32"
$ws.Cells.Item(10,2).Value = "Solana USD"
$ws.Cells.Item(10,3).Value = "SOL-USD"
$ws.Cells.Item(11,1).Value = "9.1351"
$ws.Cells.Item(11,2).Value = "Polkadot USD"
$ws.Cells.Item(11,3).Value = "DOT-USD"
$ws.Cells.Item(12,1).Value = "0.071363"
$ws.Cells.Item(12,2).Value = "Dogecoin USD"
$ws.Cells.Item(12,3).Value = "DOGE-USD"
$ws.Cells.Item(13,1).Value = "0.053941"
$ws.Cells.Item(13,2).Value = "HEX USD"
$ws.Cells.Item(13,3).Value = "HEX-USD"
$ws.Cells.Item(14,1).Value = "28.35"
$ws.Cells.Item(14,2).Value = "Avalanche USD"
$ws.Cells.Item(14,3).Value = "AVAX-USD"
$ws.Cells.Item(15,1).Value = "0.999564"
$ws.Cells.Item(15,2).Value = "Dai USD"
$ws.Cells.Item(15,3).Value = "DAI-USD"
$ws.Cells.Item(16,1).Value = "0.923054"
$ws.Cells.Item(16,2).Value = "Polygon USD"
$ws.Cells.Item(16,3).Value = "MATIC-USD"
$ws.Cells.Item(17,1).Value = "0.070251"
$ws.Cells.Item(17,2).Value = "Wrapped TRON USD"
$ws.Cells.Item(17,3).Value = "WTRX-USD"
$ws.Cells.Item(18,1).Value = "0.000012"
$ws.Cells.Item(18,2).Value = "SHIBA INU USD"
$ws.Cells.Item(18,3).Value = "SHIB-USD"
$ws.Cells.Item(19,1).Value = "85.84"
$ws.Cells.Item(19,2).Value = "STATERA USD"
$ws.Cells.Item(19,3).Value = "STA-USD"
$ws.Cells.Item(20,1).Value = "8.7848"
$ws.Cells.Item(20,2).Value = "Uniswap USD"
$ws.Cells.Item(20,3).Value = "UNI1-USD"
$ws.Cells.Item(21,1).Value = "0.070231"
$ws.Cells.Item(21,2).Value = "TRON USD"
$ws.Cells.Item(21,3).Value = "TRX-USD"
$ws.Cells.Item(22,1).Value = "1,709.10"
$ws.Cells.Item(22,2).Value = "Lido stETH USD"
$ws.Cells.Item(22,3).Value = "STETH-USD"
$ws.Cells.Item(23,1).Value = "23,799.38"
$ws.Cells.Item(23,2).Value = "Wrapped Bitcoin USD"
$ws.Cells.Item(23,3).Value = "WBTC-USD"
$ws.Cells.Item(24,1).Value = "37.75"
$ws.Cells.Item(24,2).Value = "Ethereum Classic USD"
$ws.Cells.Item(24,3).Value = "ETC-USD"
$ws.Cells.Item(25,1).Value = "4.7958"
$ws.Cells.Item(25,2).Value = "UNUS SED LEO USD"
$ws.Cells.Item(25,3).Value = "LEO-USD"
$ws.Cells.Item(26,1).Value = "62.43"
$ws.Cells.Item(26,2).Value = "Litecoin USD"
$ws.Cells.Item(26,3).Value = "LTC-USD"

# Reset column A style back to Normal (removes the temporary Text numberformat style index)
$ws.Range("A2:A26").Style = "Normal"
